$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as plain text (e.g. "28.509.10" uses
# dots as thousands separators and is not a real number), so force each
# touched Price cell to Text format before writing, to avoid Excel silently
# re-interpreting the string as a number/date and losing formatting such as
# trailing zeros (e.g. "0.5160" -> 0.516).

# Rows 33/34: Hedera and Filecoin swap places (with updated price/volume)
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.746"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07603"
$ws.Range("E34").Value = "  +7.87%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.516.72"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.09"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.56"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5160"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3873"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08324"
$ws.Range("E9").Value = "  +8.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.122"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.02"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.423"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.502"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.823.93"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.17"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001123"
$ws.Range("E18").Value = "  +3.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06658"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.82"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.057"
$ws.Range("E22").Value = "  -1.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.555.37"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.280"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.78"
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.035.54"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.415"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.26"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  +0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.099"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.684"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2237"
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.295"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.86"
$ws.Range("E39").Value = "  +5.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.771"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6405"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.194"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.61"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6134"
$ws.Range("E45").Value = "  +3.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.807"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "127.48"
$ws.Range("E47").Value = "  +2.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.003"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.207"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06987"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").Value = "  +0.49%  "
